# Auto-generated script applying cryptos.xlsx price/volume update
# commit: "Updated cryptos list on Sat Dec 30 14:46:11 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows that swapped rank position (B/C/D/E all change) ---

# Row 21: was PancakeSwap -> now Litecoin
$ws.Range("B21").Value = 'Litecoin'
$ws.Range("C21").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D21").Value = '73.31'
$ws.Range("E21").Value = '  -4.35%  '

# Row 22: was Litecoin -> now PancakeSwap
$ws.Range("B22").Value = 'PancakeSwap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D22").Value = "'3.60"
$ws.Range("E22").Value = '  -0.59%  '

# Row 27: was Toncoin -> now Cosmos
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '10.84'
$ws.Range("E27").Value = '  -6.15%  '

# Row 28: was Cosmos -> now Toncoin
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '2.43'
$ws.Range("E28").Value = '  +7.86%  '

# --- Remaining rows: price (D) and/or volume (E) updates only ---

$ws.Range("D2").Value = '42.123.36'
$ws.Range("E2").Value = '  -2.35%  '

$ws.Range("D3").Value = '2.294.66'
$ws.Range("E3").Value = '  -3.56%  '

$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").Value = '317.35'
$ws.Range("E5").Value = '  +0.02%  '

$ws.Range("D6").Value = '104.27'
$ws.Range("E6").Value = '  -4.79%  '

$ws.Range("D7").Value = '0.631'
$ws.Range("E7").Value = '  -1.32%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  -2.92%  '

$ws.Range("D10").Value = '39.52'
$ws.Range("E10").Value = '  -4.86%  '

$ws.Range("E11").Value = '  -2.66%  '

$ws.Range("E12").Value = '  -4.03%  '

$ws.Range("E13").Value = '  -0.51%  '

$ws.Range("D14").Value = '0.964'
$ws.Range("E14").Value = '  -5.18%  '

$ws.Range("E15").Value = '  -4.86%  '

$ws.Range("D16").Value = '2.640.15'
$ws.Range("E16").Value = '  -3.36%  '

$ws.Range("D17").Value = '2.295.77'
$ws.Range("E17").Value = '  -2.93%  '

$ws.Range("D18").Value = '42.099.50'
$ws.Range("E18").Value = '  -2.39%  '

$ws.Range("E19").Value = '  -5.17%  '

$ws.Range("D20").Value = '0.0000105'
$ws.Range("E20").Value = '  -1.63%  '

$ws.Range("D23").Value = '277.69'
$ws.Range("E23").Value = '  +3.13%  '

$ws.Range("D24").Value = '10.63'
$ws.Range("E24").Value = '  +10.55%  '

$ws.Range("E25").Value = '  -3.63%  '

$ws.Range("E26").Value = '  +0.36%  '

$ws.Range("D29").Value = "'22.90"
$ws.Range("E29").Value = '  -3.02%  '

$ws.Range("D30").Value = '35.89'
$ws.Range("E30").Value = '  -4.01%  '

$ws.Range("D31").Value = '163.41'
$ws.Range("E31").Value = '  -3.43%  '

$ws.Range("E32").Value = '  -4.42%  '

$ws.Range("D33").Value = '5.83'
$ws.Range("E33").Value = '  -6.18%  '

$ws.Range("D34").Value = '2.77'
$ws.Range("E34").Value = '  -5.00%  '

$ws.Range("E35").Value = '  +3.37%  '

$ws.Range("D36").Value = '0.113'
$ws.Range("E36").Value = '  -7.26%  '

$ws.Range("D37").Value = '4.58'
$ws.Range("E37").Value = '  -3.57%  '

$ws.Range("E38").Value = '  -4.60%  '

$ws.Range("D39").Value = '3.74'
$ws.Range("E39").Value = '  -3.80%  '

$ws.Range("E40").Value = '  +2.48%  '

$ws.Range("E41").Value = '  -4.23%  '

$ws.Range("D42").Value = '1.46'
$ws.Range("E42").Value = '  -5.16%  '

$ws.Range("E43").Value = '  -3.78%  '

$ws.Range("E44").Value = '  -6.12%  '

$ws.Range("E45").Value = '  -0.20%  '

$ws.Range("D46").Value = '12.02'
$ws.Range("E46").Value = '  -6.16%  '

$ws.Range("D47").Value = '111.85'
$ws.Range("E47").Value = '  -2.88%  '

$ws.Range("D48").Value = "'77.10"
$ws.Range("E48").Value = '  -4.37%  '

$ws.Range("D49").Value = '8.91'
$ws.Range("E49").Value = '  -3.24%  '

$ws.Range("D50").Value = '5.28'
$ws.Range("E50").Value = '  -5.87%  '

$ws.Range("D51").Value = '1.603.59'
$ws.Range("E51").Value = '  +1.07%  '

